# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# columns (D and G) for the first data row on both the zh-cn and de-de report
# sheets, reflecting the latest handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D2").Value = "2016-01-29 02:27:54"
$zhcn.Range("G2").Value = "2016-01-29 02:28:40"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D2").Value = "2016-01-29 02:28:06"
$dede.Range("G2").Value = "2016-01-29 02:29:00"
